# Session-01 deck rebuild: layout corrections on Slide 16
# (checklist rows get slightly taller rounded-rect backgrounds / checkbox
#  chips, and the checkbox + label rows shift down to stay vertically
#  centered in the taller rows; the closing line shifts down to match).
#
# NOTE on the literal point values below: Shape.Top/.Left/.Width/.Height
# and Adjustments.Item(n) are expressed in points in the PowerPoint object
# model (1 pt = 12700 EMU) and are stored internally as single-precision
# floats, then truncated back to EMU on save. A "clean" pt = emu/12700
# literal can therefore round-trip one EMU short of the intended target.
# The constants here are the smallest point values that truncate back to
# the exact target EMU from the diff, so the saved OOXML matches exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# --- Row 1 ("Navigated to a GitHub repository") ---------------------------
$bg = $s.Shapes.Item("Text 4")
$bg.Height = 37.99220472440945          # 482501 EMU
$bg.Adjustments.Item(1) = 0.13161        # adj val 13161

$chip = $s.Shapes.Item("Text 5")
$chip.Top = 164.23827771653544           # 2085826 EMU
$chip.Height = 25.992214724409447        # 330101 EMU
$chip.Adjustments.Item(1) = 2.77006      # adj val 277006

$check = $s.Shapes.Item("Text 6")
$check.Top = 165.0820572440945           # 2096542 EMU

$label = $s.Shapes.Item("Text 7")
$label.Top = 165.0820572440945           # 2096542 EMU

# --- Row 2 ("Read the README to understand its purpose") ------------------
$bg = $s.Shapes.Item("Text 8")
$bg.Top = 202.23048244094488             # 2568327 EMU
$bg.Height = 37.99220472440945           # 482501 EMU
$bg.Adjustments.Item(1) = 0.13161        # adj val 13161

$chip = $s.Shapes.Item("Text 9")
$chip.Top = 208.23048244094488           # 2644527 EMU
$chip.Height = 25.992214724409447        # 330101 EMU
$chip.Adjustments.Item(1) = 2.77006      # adj val 277006

$check = $s.Shapes.Item("Text 10")
$check.Top = 209.07426196850395          # 2655243 EMU

$label = $s.Shapes.Item("Text 11")
$label.Top = 209.07426196850395          # 2655243 EMU

# --- Row 3 ("Found a folder and opened it") --------------------------------
$bg = $s.Shapes.Item("Text 12")
$bg.Top = 246.22268716535433             # 3127028 EMU
$bg.Height = 37.99220472440945           # 482501 EMU
$bg.Adjustments.Item(1) = 0.13161        # adj val 13161

$chip = $s.Shapes.Item("Text 13")
$chip.Top = 252.22268716535433           # 3203228 EMU
$chip.Height = 25.992214724409447        # 330101 EMU
$chip.Adjustments.Item(1) = 2.77006      # adj val 277006

$check = $s.Shapes.Item("Text 14")
$check.Top = 253.06638795275592          # 3213943 EMU

$label = $s.Shapes.Item("Text 15")
$label.Top = 253.06638795275592          # 3213943 EMU

# --- Row 4 ("Read a file with real farm data") -----------------------------
$bg = $s.Shapes.Item("Text 16")
$bg.Top = 290.21489188976375             # 3685729 EMU
$bg.Height = 37.99220472440945           # 482501 EMU
$bg.Adjustments.Item(1) = 0.13161        # adj val 13161

$chip = $s.Shapes.Item("Text 17")
$chip.Top = 296.21489188976375           # 3761929 EMU
$chip.Height = 25.992214724409447        # 330101 EMU
$chip.Adjustments.Item(1) = 2.77006      # adj val 277006

$check = $s.Shapes.Item("Text 18")
$check.Top = 297.05858267716536          # 3772644 EMU

$label = $s.Shapes.Item("Text 19")
$label.Top = 297.05858267716536          # 3772644 EMU

# --- Closing line ("That's the foundation of everything we'll learn.") ----
$closing = $s.Shapes.Item("Text 20")
$closing.Top = 336.1992125984252         # 4269730 EMU
